$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The workbook tracks LMNX yearly financials (Income Statement, Balance
# Sheet, Cash Flow) with one column per fiscal year, most-recent-year
# first starting at column D. This update adds the newest fiscal year
# (period ending 2018-12-31) as a new column D, pushing the existing
# FY2017..FY2011 columns from D:K over to E:L.
# ---------------------------------------------------------------------------

# Insert a new blank column before D; existing D:K data shifts to E:L.
$ws.Columns("D").Insert()

# The freshly inserted column D has no number formatting yet - copy the
# formatting (date format on the "Period Ending" rows, #,##0 elsewhere)
# from column E (the old column D, now shifted one to the right).
$ws.Columns("E").Copy()
$ws.Columns("D").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Populate new column D with FY2018 figures for each statement line.
# ---------------------------------------------------------------------------

# --- Income Statement (rows 7-35) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 315800
$ws.Range("D9").Value = 120300
$ws.Range("D10").Value = 195500
$ws.Range("D12").Value = 47200
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 8700
$ws.Range("D17").Value = 288000
$ws.Range("D18").Value = 27800
$ws.Range("D20").Value = 500
$ws.Range("D21").Value = 52000
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = 28300
$ws.Range("D24").Value = 16500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 11800
$ws.Range("D27").Value = 11500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 6700
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -500
$ws.Range("D33").Value = 18200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 18200

# --- Balance Sheet (rows 38-77) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 76400
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 53400
$ws.Range("D44").Value = 63300
$ws.Range("D45").Value = 9700
$ws.Range("D46").Value = 202700
$ws.Range("D47").Value = 2800
$ws.Range("D48").Value = 66300
$ws.Range("D49").Value = 229900
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 23400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 525200
$ws.Range("D57").Value = 14500
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 36900
$ws.Range("D60").Value = 51400
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 6100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 57500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 103400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 467700
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (rows 80-102) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 18200
$ws.Range("D83").Value = 23700
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 50900
$ws.Range("D91").Value = -21300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -93500
$ws.Range("D96").Value = -10700
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -8400
$ws.Range("D101").Value = 300
$ws.Range("D102").Value = -50700

# Row 91 ("Capital Expenditures") was also restated for the prior years at
# the same time - not a pure shift, so fix up E:K explicitly.
$ws.Range("E91").Value = -14600
$ws.Range("F91").Value = -13100
$ws.Range("G91").Value = -18700
$ws.Range("H91").Value = -17100
$ws.Range("I91").Value = -18100
$ws.Range("J91").Value = -9800
$ws.Range("K91").Value = -11400
